$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows (4-136) down to (5-137)
$ws.Rows("4").Insert()

# Copy the (now shifted) old row 4 data - which is now in row 5 - into the new row 4
$ws.Range("A5:R5").Copy()
$ws.Range("A4:R4").PasteSpecial()

# Set the new date for row 4
$ws.Range("D4").Value = 44643

Write-Host "done"
